$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("data2 weibull")
$ws.Range("E2").Value = -2.33006830871177
$ws.Range("F2").Value = 0.10533870155035
$ws.Range("G2").Value = 0.105924771630349
$ws.Range("H2").Value = 0.0881247958125287
$ws.Range("I2").Value = 0.0110962420443137
$ws.Range("J2").Value = 0.00776597963699988
$ws.Range("K2").Value = -0.00350576567788289
$ws.Range("E3").Value = -2.49972872840229
$ws.Range("F3").Value = 0.122057447093382
$ws.Range("G3").Value = 0.0391414385130606
$ws.Range("H3").Value = 0.0517908396172799
$ws.Range("I3").Value = 0.0148980203909537
$ws.Range("J3").Value = 0.00268229106826281
$ws.Range("K3").Value = -0.00257155016007909
$ws.Range("E4").Value = -3.09126802584416
$ws.Range("F4").Value = 0.175312283413101
$ws.Range("G4").Value = 0.143473106550909
$ws.Range("H4").Value = 0.117818399725227
$ws.Range("I4").Value = 0.0307343967155154
$ws.Range("J4").Value = 0.0138811753138134
$ws.Range("K4").Value = -0.0141763938303875
$ws.Range("E5").Value = -3.48734749630726
$ws.Range("F5").Value = 0.242367660737922
$ws.Range("G5").Value = 0.282416396265322
$ws.Range("H5").Value = 0.11231437447069
$ws.Range("I5").Value = 0.0587420829715727
$ws.Range("J5").Value = 0.0126145187127423
$ws.Range("K5").Value = -0.0229317679615367
$ws.Range("E6").Value = -2.41945335488241
$ws.Range("F6").Value = 0.107605083706125
$ws.Range("G6").Value = -0.0399294119806398
$ws.Range("H6").Value = 0.0802228258285359
$ws.Range("I6").Value = 0.0115788540394022
$ws.Range("J6").Value = 0.00643570178391561
$ws.Range("K6").Value = -0.00318512887970389
$ws.Range("E7").Value = -2.56891849151363
$ws.Range("F7").Value = 0.0953207574839895
$ws.Range("G7").Value = -0.0434775808987697
$ws.Range("H7").Value = 0.101936547928577
$ws.Range("I7").Value = 0.00908604680732155
$ws.Range("J7").Value = 0.0103910598035951
$ws.Range("K7").Value = -0.00562530671538125
$ws.Range("E8").Value = -2.28927585322432
$ws.Range("F8").Value = 0.263546347264198
$ws.Range("G8").Value = -0.107253607990594
$ws.Range("H8").Value = 0.169984735668374
$ws.Range("I8").Value = 0.0694566771563014
$ws.Range("J8").Value = 0.0288948103602471
$ws.Range("K8").Value = -0.0143544642403251
$ws.Range("E9").Value = -2.97605246674702
$ws.Range("F9").Value = 0.323540867216825
$ws.Range("G9").Value = 0.137278719565766
$ws.Range("H9").Value = 0.216077907045135
$ws.Range("I9").Value = 0.104678692759415
$ws.Range("J9").Value = 0.0466896619130061
$ws.Range("K9").Value = -0.0591706299533721
$ws.Range("E10").Value = -2.0472551672577
$ws.Range("F10").Value = 0.261938085818681
$ws.Range("G10").Value = 0.155430830791449
$ws.Range("H10").Value = 0.172266877897778
$ws.Range("I10").Value = 0.0686115608023547
$ws.Range("J10").Value = 0.0296758772206481
$ws.Range("K10").Value = -0.0339121655313017
$ws.Range("E11").Value = -2.6879308622404
$ws.Range("F11").Value = 0.254893879407089
$ws.Range("G11").Value = 0.130167622053188
$ws.Range("H11").Value = 0.11158061989941
$ws.Range("I11").Value = 0.0649708897591956
$ws.Range("J11").Value = 0.0124502347371367
$ws.Range("K11").Value = -0.0206858569923194
$ws.Range("E12").Value = -2.76488446758512
$ws.Range("F12").Value = 0.2959383108534
$ws.Range("G12").Value = 0.174964827340232
$ws.Range("H12").Value = 0.114840037842542
$ws.Range("I12").Value = 0.0875794838307637
$ws.Range("J12").Value = 0.0131882342916765
$ws.Range("K12").Value = -0.0198566769345156
$ws.Range("E13").Value = -3.92419203737332
$ws.Range("F13").Value = 0.494772230509085
$ws.Range("G13").Value = 0.360426783558681
$ws.Range("H13").Value = 0.187127507302931
$ws.Range("I13").Value = 0.244799560082935
$ws.Range("J13").Value = 0.0350167039894086
$ws.Range("K13").Value = -0.0810677652745029
$ws.Range("E14").Value = -2.79738339481437
$ws.Range("F14").Value = 0.171698196289065
$ws.Range("G14").Value = -0.0139801110920569
$ws.Range("H14").Value = 0.105889677952872
$ws.Range("I14").Value = 0.0294802706089184
$ws.Range("J14").Value = 0.011212623896963
$ws.Range("K14").Value = -0.0116199924086907
$ws.Range("E15").Value = -2.95410791709162
$ws.Range("F15").Value = 0.173133254677496
$ws.Range("G15").Value = 0.0556593266224006
$ws.Range("H15").Value = 0.0926145582214106
$ws.Range("I15").Value = 0.0299751238752225
$ws.Range("J15").Value = 0.00857745639454706
$ws.Range("K15").Value = -0.0112536040524681
$ws.Range("E16").Value = -2.79756053367509
$ws.Range("F16").Value = 0.119642980975563
$ws.Range("G16").Value = -0.0884476008590871
$ws.Range("H16").Value = 0.0838583312620675
$ws.Range("I16").Value = 0.014314442896719
$ws.Range("J16").Value = 0.00703221972205866
$ws.Range("K16").Value = -0.00234476540435095
$ws.Range("E17").Value = -1.8496768167544
$ws.Range("F17").Value = 0.0713815590614557
$ws.Range("G17").Value = -0.201232123211273
$ws.Range("H17").Value = 0.036854147082554
$ws.Range("I17").Value = 0.00509532697404409
$ws.Range("J17").Value = 0.00135822815718253
$ws.Range("K17").Value = -0.000288346995536712
$ws.Range("E18").Value = -2.04690973252514
$ws.Range("F18").Value = 0.0737316968516354
$ws.Range("G18").Value = -0.171946395700535
$ws.Range("H18").Value = 0.0388523711236188
$ws.Range("I18").Value = 0.00543636312062147
$ws.Range("J18").Value = 0.00150950674192741
$ws.Range("K18").Value = -0.00158773213442984
$ws.Range("E19").Value = -2.83545108599723
$ws.Range("F19").Value = 0.159160874987795
$ws.Range("G19").Value = 0.145673114301433
$ws.Range("H19").Value = 0.0977476357725525
$ws.Range("I19").Value = 0.0253321841268806
$ws.Range("J19").Value = 0.00955460029912358
$ws.Range("K19").Value = -0.0120621687812353
$ws.Range("E20").Value = -2.72194870874295
$ws.Range("F20").Value = 0.163891712200428
$ws.Range("G20").Value = 0.255836150826446
$ws.Range("H20").Value = 0.118697956862633
$ws.Range("I20").Value = 0.026860493327988
$ws.Range("J20").Value = 0.0140892049633635
$ws.Range("K20").Value = -0.0149839404293542
$ws.Range("E21").Value = -1.91600227134157
$ws.Range("F21").Value = 0.257928758189887
$ws.Range("G21").Value = -0.139560850406393
$ws.Range("H21").Value = 0.105925170492656
$ws.Range("I21").Value = 0.0665272443013772
$ws.Range("J21").Value = 0.0112201417438982
$ws.Range("K21").Value = -0.0112609082153799
$ws.Range("E22").Value = -2.61515646363081
$ws.Range("F22").Value = 0.320756770614934
$ws.Range("G22").Value = 0.213278475676705
$ws.Range("H22").Value = 0.172728878584829
$ws.Range("I22").Value = 0.102884905895322
$ws.Range("J22").Value = 0.0298352654971725
$ws.Range("K22").Value = -0.0348334841073863

$ws = $wb.Worksheets.Item("data2 lognormal")
$ws.Range("E2").Value = 1.97970353858253
$ws.Range("F2").Value = 0.145690465580844
$ws.Range("G2").Value = -1.12832256940951
$ws.Range("H2").Value = 0.103194604784936
$ws.Range("I2").Value = 0.0212257117611632
$ws.Range("J2").Value = 0.0106491264567192
$ws.Range("K2").Value = -0.0124958410297487
$ws.Range("E3").Value = 1.93563637872471
$ws.Range("F3").Value = 0.121436695299748
$ws.Range("G3").Value = -0.991804474991053
$ws.Range("H3").Value = 0.0583512054430059
$ws.Range("I3").Value = 0.0147468709653237
$ws.Range("J3").Value = 0.00340486317665188
$ws.Range("K3").Value = -0.00444446528904442
$ws.Range("E4").Value = 2.46833408038063
$ws.Range("F4").Value = 0.22714291112948
$ws.Range("G4").Value = -1.03846010749416
$ws.Range("H4").Value = 0.115243977650148
$ws.Range("I4").Value = 0.0515939020763749
$ws.Range("J4").Value = 0.0132811743846279
$ws.Range("K4").Value = -0.0232637853881486
$ws.Range("E5").Value = 2.81382207634397
$ws.Range("F5").Value = 0.242498989040528
$ws.Range("G5").Value = -1.1205857261163
$ws.Range("H5").Value = 0.101912722996449
$ws.Range("I5").Value = 0.0588057596856783
$ws.Range("J5").Value = 0.010386203108551
$ws.Range("K5").Value = -0.021909225552673
$ws.Range("E6").Value = 2.01849548600067
$ws.Range("F6").Value = 0.122883242611463
$ws.Range("G6").Value = -1.01997434819274
$ws.Range("H6").Value = 0.0818857483036436
$ws.Range("I6").Value = 0.0151002913147078
$ws.Range("J6").Value = 0.00670527577524767
$ws.Range("K6").Value = -0.00742396882898494
$ws.Range("E7").Value = 2.3420142874128
$ws.Range("F7").Value = 0.215810175558297
$ws.Range("G7").Value = -1.08629849598012
$ws.Range("H7").Value = 0.13211005024603
$ws.Range("I7").Value = 0.0465740318745031
$ws.Range("J7").Value = 0.0174530653760086
$ws.Range("K7").Value = -0.0269241863120591
$ws.Range("E8").Value = 1.52982703612423
$ws.Range("F8").Value = 0.343917875306103
$ws.Range("G8").Value = -0.830115526445352
$ws.Range("H8").Value = 0.186708829854108
$ws.Range("I8").Value = 0.118279504955064
$ws.Range("J8").Value = 0.0348601871454903
$ws.Range("K8").Value = -0.0477083810907798
$ws.Range("E9").Value = 2.4381676875878
$ws.Range("F9").Value = 0.484812865854487
$ws.Range("G9").Value = -1.07201620635763
$ws.Range("H9").Value = 0.192364462982496
$ws.Range("I9").Value = 0.235043514898041
$ws.Range("J9").Value = 0.0370040866185441
$ws.Range("K9").Value = -0.0878022120874772
$ws.Range("E10").Value = 1.41853795862842
$ws.Range("F10").Value = 0.298670614418002
$ws.Range("G10").Value = -0.997272275177952
$ws.Range("H10").Value = 0.129046366793532
$ws.Range("I10").Value = 0.089204135916827
$ws.Range("J10").Value = 0.0166529647826108
$ws.Range("K10").Value = -0.0321888868121076
$ws.Range("E11").Value = 1.91397156346093
$ws.Range("F11").Value = 0.251259573281716
$ws.Range("G11").Value = -0.965044938060798
$ws.Range("H11").Value = 0.0796827302863566
$ws.Range("I11").Value = 0.0631313731657098
$ws.Range("J11").Value = 0.00634933750588825
$ws.Range("K11").Value = -0.0158753027847061
$ws.Range("E12").Value = 1.98537693246169
$ws.Range("F12").Value = 0.309598793339315
$ws.Range("G12").Value = -0.971956116966394
$ws.Range("H12").Value = 0.12666401900976
$ws.Range("I12").Value = 0.0958514128371598
$ws.Range("J12").Value = 0.0160437737117049
$ws.Range("K12").Value = -0.0289126257453649
$ws.Range("E13").Value = 2.81513530583061
$ws.Range("F13").Value = 0.549529613637574
$ws.Range("G13").Value = -1.03564365282378
$ws.Range("H13").Value = 0.177278779745409
$ws.Range("I13").Value = 0.301982796264662
$ws.Range("J13").Value = 0.0314277657480212
$ws.Range("K13").Value = -0.0905067219448786
$ws.Range("E14").Value = 2.20329872621066
$ws.Range("F14").Value = 0.251472192731928
$ws.Range("G14").Value = -0.949589611386702
$ws.Range("H14").Value = 0.120483913853294
$ws.Range("I14").Value = 0.0632382637174039
$ws.Range("J14").Value = 0.0145163734974079
$ws.Range("K14").Value = -0.0269335727041458
$ws.Range("E15").Value = 2.40219799434963
$ws.Range("F15").Value = 0.234519749464493
$ws.Range("G15").Value = -1.00901173307717
$ws.Range("H15").Value = 0.102287905033608
$ws.Range("I15").Value = 0.0549995128888884
$ws.Range("J15").Value = 0.0104628155161643
$ws.Range("K15").Value = -0.0211596960556413
$ws.Range("E16").Value = 2.20170535405669
$ws.Range("F16").Value = 0.21892158455339
$ws.Range("G16").Value = -0.902310598510134
$ws.Range("H16").Value = 0.114508545878668
$ws.Range("I16").Value = 0.047926660183367
$ws.Range("J16").Value = 0.013112207079247
$ws.Range("K16").Value = -0.0217737793079823
$ws.Range("E17").Value = 0.999415283650382
$ws.Range("F17").Value = 0.0710460548893685
$ws.Range("G17").Value = -0.710160646468684
$ws.Range("H17").Value = 0.0371536156696288
$ws.Range("I17").Value = 0.00504754191534316
$ws.Range("J17").Value = 0.00138039115732648
$ws.Range("K17").Value = -0.00121965787689204
$ws.Range("E18").Value = 1.13237362326564
$ws.Range("F18").Value = 0.0882260712790033
$ws.Range("G18").Value = -0.691559258068776
$ws.Range("H18").Value = 0.0375498739940975
$ws.Range("I18").Value = 0.00778383965332777
$ws.Range("J18").Value = 0.0014099930369726
$ws.Range("K18").Value = -0.00251196077708807
$ws.Range("E19").Value = 2.25372991308615
$ws.Range("F19").Value = 0.193599458785533
$ws.Range("G19").Value = -1.05166003595957
$ws.Range("H19").Value = 0.0926540495239974
$ws.Range("I19").Value = 0.0374807504420513
$ws.Range("J19").Value = 0.00858477289319537
$ws.Range("K19").Value = -0.0159888332730207
$ws.Range("E20").Value = 2.29516138426015
$ws.Range("F20").Value = 0.222364136608754
$ws.Range("G20").Value = -1.17781611095422
$ws.Range("H20").Value = 0.104923760248021
$ws.Range("I20").Value = 0.0494458092497564
$ws.Range("J20").Value = 0.0110089954645841
$ws.Range("K20").Value = -0.0206429316639446
$ws.Range("E21").Value = 1.02869503297525
$ws.Range("F21").Value = 0.265769224642769
$ws.Range("G21").Value = -0.717688599529452
$ws.Range("H21").Value = 0.103533830249769
$ws.Range("I21").Value = 0.0706332807672187
$ws.Range("J21").Value = 0.010719254006188
$ws.Range("K21").Value = -0.0171998828454053
$ws.Range("E22").Value = 1.85768941713106
$ws.Range("F22").Value = 0.336969112651298
$ws.Range("G22").Value = -0.993295258875856
$ws.Range("H22").Value = 0.157840161947571
$ws.Range("I22").Value = 0.113548182881003
$ws.Range("J22").Value = 0.0249135167236354
$ws.Range("K22").Value = -0.0414908144906914

$ws = $wb.Worksheets.Item("data2 llogis")
$ws.Range("E2").Value = -1.72519678835276
$ws.Range("F2").Value = 0.0869111420140341
$ws.Range("G2").Value = 1.99546287416998
$ws.Range("H2").Value = 0.204557243415721
$ws.Range("I2").Value = 0.00755354660618361
$ws.Range("J2").Value = 0.0418436658338387
$ws.Range("K2").Value = 0.00949575258569563
$ws.Range("E3").Value = -1.93060133863177
$ws.Range("F3").Value = 0.102811464864021
$ws.Range("G3").Value = 1.71540910960479
$ws.Range("H3").Value = 0.128615701649951
$ws.Range("I3").Value = 0.0105701973074859
$ws.Range("J3").Value = 0.0165419987109093
$ws.Range("K3").Value = 0.00506328286632033
$ws.Range("E4").Value = -2.3370003112681
$ws.Range("F4").Value = 0.129544185691206
$ws.Range("G4").Value = 1.73953940198937
$ws.Range("H4").Value = 0.208695357802345
$ws.Range("I4").Value = 0.0167816960463977
$ws.Range("J4").Value = 0.0435537523682489
$ws.Range("K4").Value = 0.0130353650991837
$ws.Range("E5").Value = -2.48934273642372
$ws.Range("F5").Value = 0.116141918019508
$ws.Range("G5").Value = 1.86641737773398
$ws.Range("H5").Value = 0.18441970056748
$ws.Range("I5").Value = 0.0134889451212501
$ws.Range("J5").Value = 0.034010625957399
$ws.Range("K5").Value = 0.00574248268227631
$ws.Range("E6").Value = -1.94481970743779
$ws.Range("F6").Value = 0.105783665680401
$ws.Range("G6").Value = 1.78244840612281
$ws.Range("H6").Value = 0.145920306625554
$ws.Range("I6").Value = 0.0111901839247827
$ws.Range("J6").Value = 0.0212927358856957
$ws.Range("K6").Value = 0.00854226654679157
$ws.Range("E7").Value = -2.08438385665272
$ws.Range("F7").Value = 0.0833283834066869
$ws.Range("G7").Value = 2.00821523370502
$ws.Range("H7").Value = 0.264693199286976
$ws.Range("I7").Value = 0.00694361948117181
$ws.Range("J7").Value = 0.0700624897487748
$ws.Range("K7").Value = 0.01381622202849
$ws.Range("E8").Value = -1.83920081528535
$ws.Range("F8").Value = 0.283458625754857
$ws.Range("G8").Value = 1.43323566724509
$ws.Range("H8").Value = 0.380173870798206
$ws.Range("I8").Value = 0.0803487925148324
$ws.Range("J8").Value = 0.144532172037691
$ws.Range("K8").Value = 0.0219698548481498
$ws.Range("E9").Value = -2.28588536418224
$ws.Range("F9").Value = 0.118148159030896
$ws.Range("G9").Value = 1.97834259892908
$ws.Range("H9").Value = 0.38226490490551
$ws.Range("I9").Value = 0.01395898748239
$ws.Range("J9").Value = 0.146126457522419
$ws.Range("K9").Value = -0.00159880999586965
$ws.Range("E10").Value = -1.45210266844635
$ws.Range("F10").Value = 0.163509255779092
$ws.Range("G10").Value = 1.76506582427408
$ws.Range("H10").Value = 0.238246846292879
$ws.Range("I10").Value = 0.0267352767254325
$ws.Range("J10").Value = 0.0567615597685027
$ws.Range("K10").Value = -0.0101279729218984
$ws.Range("E11").Value = -2.0339517655342
$ws.Range("F11").Value = 0.179930231478184
$ws.Range("G11").Value = 1.64617655510194
$ws.Range("H11").Value = 0.157352972748654
$ws.Range("I11").Value = 0.0323748881997928
$ws.Range("J11").Value = 0.0247599580328386
$ws.Range("K11").Value = -0.00937185231836174
$ws.Range("E12").Value = -2.07102831617382
$ws.Range("F12").Value = 0.21136996067186
$ws.Range("G12").Value = 1.68163101716978
$ws.Range("H12").Value = 0.224759741543284
$ws.Range("I12").Value = 0.0446772602744237
$ws.Range("J12").Value = 0.0505169414186039
$ws.Range("K12").Value = 0.00891023701979488
$ws.Range("E13").Value = -2.76933170304976
$ws.Range("F13").Value = 0.177956593636583
$ws.Range("G13").Value = 1.91020861268176
$ws.Range("H13").Value = 0.371607190811063
$ws.Range("I13").Value = 0.0316685492187359
$ws.Range("J13").Value = 0.13809190426249
$ws.Range("K13").Value = -0.0048304199862856
$ws.Range("E14").Value = -2.29189352078826
$ws.Range("F14").Value = 0.134661724775328
$ws.Range("G14").Value = 1.63959728057625
$ws.Range("H14").Value = 0.243901739439504
$ws.Range("I14").Value = 0.0181337801194661
$ws.Range("J14").Value = 0.0594880585016157
$ws.Range("K14").Value = 0.0106255584561674
$ws.Range("E15").Value = -2.35918673923445
$ws.Range("F15").Value = 0.116437894083279
$ws.Range("G15").Value = 1.75657302549737
$ws.Range("H15").Value = 0.210949146448444
$ws.Range("I15").Value = 0.013557783178549
$ws.Range("J15").Value = 0.0444995423873269
$ws.Range("K15").Value = 0.00663085164219347
$ws.Range("E16").Value = -2.38452510033337
$ws.Range("F16").Value = 0.149022890586942
$ws.Range("G16").Value = 1.58591775413186
$ws.Range("H16").Value = 0.246475210808504
$ws.Range("I16").Value = 0.0222078219188878
$ws.Range("J16").Value = 0.0607500295430965
$ws.Range("K16").Value = 0.0211604529847967
$ws.Range("E17").Value = -1.37351061111828
$ws.Range("F17").Value = 0.0975451183537697
$ws.Range("G17").Value = 1.17954533219363
$ws.Range("H17").Value = 0.064612489657263
$ws.Range("I17").Value = 0.00951505011465093
$ws.Range("J17").Value = 0.00417477381970992
$ws.Range("K17").Value = 0.00140889801365611
$ws.Range("E18").Value = -1.64016821659955
$ws.Range("F18").Value = 0.0855002121713666
$ws.Range("G18").Value = 1.1820943127845
$ws.Range("H18").Value = 0.0656643713394952
$ws.Range("I18").Value = 0.0073102862813487
$ws.Range("J18").Value = 0.00431180966341111
$ws.Range("K18").Value = -0.000502977650851661
$ws.Range("E19").Value = -2.14081736529117
$ws.Range("F19").Value = 0.0953969328983439
$ws.Range("G19").Value = 1.80317848000473
$ws.Range("H19").Value = 0.177904693938904
$ws.Range("I19").Value = 0.00910057480641112
$ws.Range("J19").Value = 0.0316500801254951
$ws.Range("K19").Value = 0.00299024922680052
$ws.Range("E20").Value = -1.97028949193194
$ws.Range("F20").Value = 0.0801902518794047
$ws.Range("G20").Value = 2.12084194642728
$ws.Range("H20").Value = 0.213035641702335
$ws.Range("I20").Value = 0.00643047649648237
$ws.Range("J20").Value = 0.0453841846355257
$ws.Range("K20").Value = -0.00257506918792403
$ws.Range("E21").Value = -1.46716693033476
$ws.Range("F21").Value = 0.287705749115125
$ws.Range("G21").Value = 1.26059972867232
$ws.Range("H21").Value = 0.170502953872554
$ws.Range("I21").Value = 0.0827745980738954
$ws.Range("J21").Value = 0.0290712572792663
$ws.Range("K21").Value = 0.0044798074636298
$ws.Range("E22").Value = -1.89524984654235
$ws.Range("F22").Value = 0.220566185666429
$ws.Range("G22").Value = 1.75319588941478
$ws.Range("H22").Value = 0.328357094076466
$ws.Range("I22").Value = 0.0486494422594377
$ws.Range("J22").Value = 0.107818381230341
$ws.Range("K22").Value = 0.00702714895855935

$ws = $wb.Worksheets.Item("data2 gompertz")
$ws.Range("E2").Value = -2.01765690253255
$ws.Range("F2").Value = 0.103447449067094
$ws.Range("G2").Value = -0.0190740131440861
$ws.Range("H2").Value = 0.0138713780231997
$ws.Range("I2").Value = 0.0107013747184891
$ws.Range("J2").Value = 0.000192415128262507
$ws.Range("K2").Value = -0.00050452431854454
$ws.Range("E3").Value = -2.29339372378354
$ws.Range("F3").Value = 0.132301106530333
$ws.Range("G3").Value = -0.013784129933345
$ws.Range("H3").Value = 0.00673971626319614
$ws.Range("I3").Value = 0.0175035827891505
$ws.Range("J3").Value = 0.0000454237753083905
$ws.Range("K3").Value = -0.00054132285507576
$ws.Range("E4").Value = -2.66618220861419
$ws.Range("F4").Value = 0.144808633201942
$ws.Range("G4").Value = -0.0198690908444268
$ws.Range("H4").Value = 0.0191775399753587
$ws.Range("I4").Value = 0.0209695402498145
$ws.Range("J4").Value = 0.00036777803950648
$ws.Range("K4").Value = -0.00167073607669963
$ws.Range("E5").Value = -3.00530460641232
$ws.Range("F5").Value = 0.175142582567864
$ws.Range("G5").Value = 0.00745643115831869
$ws.Range("H5").Value = 0.0142202203186205
$ws.Range("I5").Value = 0.0306749242285411
$ws.Range("J5").Value = 0.000202214665910108
$ws.Range("K5").Value = -0.00183131118844048
$ws.Range("E6").Value = -2.17249190145231
$ws.Range("F6").Value = 0.111546443989113
$ws.Range("G6").Value = -0.0268888286602143
$ws.Range("H6").Value = 0.00923898580002034
$ws.Range("I6").Value = 0.0124426091666163
$ws.Range("J6").Value = 0.0000853588586129775
$ws.Range("K6").Value = -0.000441079631820998
$ws.Range("E7").Value = -2.25929198323504
$ws.Range("F7").Value = 0.08253900642319
$ws.Range("G7").Value = -0.0299328683024798
$ws.Range("H7").Value = 0.0095122801950541
$ws.Range("I7").Value = 0.0068126875813274
$ws.Range("J7").Value = 0.0000904834745092185
$ws.Range("K7").Value = -0.000146717072827743
$ws.Range("E8").Value = -2.15271197631863
$ws.Range("F8").Value = 0.295240845282059
$ws.Range("G8").Value = -0.0312207305439134
$ws.Range("H8").Value = 0.0287361000618847
$ws.Range("I8").Value = 0.0871671567228649
$ws.Range("J8").Value = 0.000825763446766651
$ws.Range("K8").Value = -0.00352109315215479
$ws.Range("E9").Value = -2.65027185824501
$ws.Range("F9").Value = 0.225090686544683
$ws.Range("G9").Value = 0.00000447698641238702
$ws.Range("H9").Value = 0.0270325818460841
$ws.Range("I9").Value = 0.0506658171691567
$ws.Range("J9").Value = 0.000730760481265237
$ws.Range("K9").Value = -0.00415997704233024
$ws.Range("E10").Value = -1.9226526959951
$ws.Range("F10").Value = 0.2434248993688
$ws.Range("G10").Value = 0.0148059705014051
$ws.Range("H10").Value = 0.0463504303092982
$ws.Range("I10").Value = 0.0592556816327103
$ws.Range("J10").Value = 0.00214836238985711
$ws.Range("K10").Value = -0.00869296177512084
$ws.Range("E11").Value = -2.58753790062641
$ws.Range("F11").Value = 0.210820747593303
$ws.Range("G11").Value = 0.0146205057425297
$ws.Range("H11").Value = 0.0191134696191717
$ws.Range("I11").Value = 0.0444453876157992
$ws.Range("J11").Value = 0.000365324720882999
$ws.Range("K11").Value = -0.00225609370786197
$ws.Range("E12").Value = -2.54933555186575
$ws.Range("F12").Value = 0.249008871886208
$ws.Range("G12").Value = 0.0102466828139336
$ws.Range("H12").Value = 0.0195441948683702
$ws.Range("I12").Value = 0.0620054182780418
$ws.Range("J12").Value = 0.000381975553052829
$ws.Range("K12").Value = -0.00214364487714964
$ws.Range("E13").Value = -3.42672233937081
$ws.Range("F13").Value = 0.243199656921414
$ws.Range("G13").Value = 0.0243636821714332
$ws.Range("H13").Value = 0.0122956675768103
$ws.Range("I13").Value = 0.0591460731266934
$ws.Range("J13").Value = 0.000151183441159423
$ws.Range("K13").Value = -0.00126953734763064
$ws.Range("E14").Value = -2.49902757553759
$ws.Range("F14").Value = 0.151796300191924
$ws.Range("G14").Value = -0.0260432363531547
$ws.Range("H14").Value = 0.012122266923568
$ws.Range("I14").Value = 0.0230421167519567
$ws.Range("J14").Value = 0.000146949355366231
$ws.Range("K14").Value = -0.000973927502218067
$ws.Range("E15").Value = -2.60453622180521
$ws.Range("F15").Value = 0.140775675694235
$ws.Range("G15").Value = -0.0179874448132232
$ws.Range("H15").Value = 0.00981579720341566
$ws.Range("I15").Value = 0.0198177908671686
$ws.Range("J15").Value = 0.0000963498747385826
$ws.Range("K15").Value = -0.000793480139508984
$ws.Range("E16").Value = -2.47257403697872
$ws.Range("F16").Value = 0.1067374594493
$ws.Range("G16").Value = -0.0386901786342452
$ws.Range("H16").Value = 0.00949365882628895
$ws.Range("I16").Value = 0.011392885249691
$ws.Range("J16").Value = 0.000090129557909974
$ws.Range("K16").Value = -0.0000163124782557576
$ws.Range("E17").Value = -1.59135908631409
$ws.Range("F17").Value = 0.0932982144745055
$ws.Range("G17").Value = -0.0889611671755879
$ws.Range("H17").Value = 0.0166836023152949
$ws.Range("I17").Value = 0.00870455682413083
$ws.Range("J17").Value = 0.000278342586214914
$ws.Range("K17").Value = -0.00105036762195595
$ws.Range("E18").Value = -1.86178907503081
$ws.Range("F18").Value = 0.0914459151885326
$ws.Range("G18").Value = -0.0646415225040878
$ws.Range("H18").Value = 0.0128324640774315
$ws.Range("I18").Value = 0.0083623554046683
$ws.Range("J18").Value = 0.000164672134298569
$ws.Range("K18").Value = -0.00089930462489449
$ws.Range("E19").Value = -2.53105087692342
$ws.Range("F19").Value = 0.130124689890718
$ws.Range("G19").Value = -0.00238000384600223
$ws.Range("H19").Value = 0.0140000910362655
$ws.Range("I19").Value = 0.0169324349191557
$ws.Range("J19").Value = 0.000196002549023721
$ws.Range("K19").Value = -0.00125060303915194
$ws.Range("E20").Value = -2.29936955689108
$ws.Range("F20").Value = 0.105835268383319
$ws.Range("G20").Value = 0.00178324955458855
$ws.Range("H20").Value = 0.0170697952170677
$ws.Range("I20").Value = 0.0112011040337692
$ws.Range("J20").Value = 0.000291377908752626
$ws.Range("K20").Value = -0.0010767979983474
$ws.Range("E21").Value = -1.77351822022407
$ws.Range("F21").Value = 0.320554106187069
$ws.Range("G21").Value = -0.0531899742245637
$ws.Range("H21").Value = 0.0449845807907303
$ws.Range("I21").Value = 0.10275493499339
$ws.Range("J21").Value = 0.00202361250891774
$ws.Range("K21").Value = -0.0106890316796556
$ws.Range("E22").Value = -2.39640552258822
$ws.Range("F22").Value = 0.274374733260961
$ws.Range("G22").Value = 0.0140410826723287
$ws.Range("H22").Value = 0.0322726576179848
$ws.Range("I22").Value = 0.0752814942520235
$ws.Range("J22").Value = 0.00104152442972767
$ws.Range("K22").Value = -0.00475935260339732
